# Week 15 simulations added.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append this week's simulated yardage samples to the 4 running
# lists (OFF Rush, OFF Pass, DEF Rush, DEF Pass).
# ---------------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value2 = $ydsWs.Range("B2").Value2 + " 2 1 3 3 2 5 1 6 15 1 2 13 -2 5 6 3 4 11 8 11 1 2 7 1 0 -2 1"
$ydsWs.Range("B3").Value2 = $ydsWs.Range("B3").Value2 + " -1 8 3 9 4 34 5 5 9 7 6 27 9 1 12 5 7 13 17 3 3 4"
$ydsWs.Range("C2").Value2 = $ydsWs.Range("C2").Value2 + " 5 1 4 5 4 3 -4 -1 7 2 18 4 8 5 0 11 13 21 5 3 2 1 1"
$ydsWs.Range("C3").Value2 = $ydsWs.Range("C3").Value2 + " 11 2 -1 5 3 5 4 9 4 20 10 32 4 8 -3 20 36 12 3 11 2 6 5 6 4 8 14"

# ---------------------------------------------------------------------------
# OFF sheet: updated running totals after Week 15.
# ---------------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value2 = 446
$offWs.Range("D2").Value2 = 21
$offWs.Range("E2").Value2 = 21
$offWs.Range("F2").Value2 = 123
$offWs.Range("G2").Value2 = 98
$offWs.Range("H2").Value2 = 13
$offWs.Range("I2").Value2 = 14
$offWs.Range("J2").Value2 = 65
$offWs.Range("L2").Value2 = 438
$offWs.Range("M2").Value2 = 276
$offWs.Range("Q2").Value2 = 948

$offWs.Range("B3").Value2 = 21
$offWs.Range("C3").Value2 = 270
$offWs.Range("D3").Value2 = 13
$offWs.Range("F3").Value2 = 198
$offWs.Range("G3").Value2 = 63
$offWs.Range("H3").Value2 = 52
$offWs.Range("I3").Value2 = 89
$offWs.Range("J3").Value2 = 88

# ---------------------------------------------------------------------------
# DEF sheet: updated running totals after Week 15.
# ---------------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("B2").Value2 = 13
$defWs.Range("C2").Value2 = 337
$defWs.Range("D2").Value2 = 29
$defWs.Range("E2").Value2 = 21
$defWs.Range("F2").Value2 = 109
$defWs.Range("G2").Value2 = 99
$defWs.Range("I2").Value2 = 16
$defWs.Range("J2").Value2 = 59
$defWs.Range("L2").Value2 = 482
$defWs.Range("M2").Value2 = 302
$defWs.Range("O2").Value2 = 51
$defWs.Range("P2").Value2 = 38
$defWs.Range("Q2").Value2 = 898

$defWs.Range("C3").Value2 = 365
$defWs.Range("E3").Value2 = 37
$defWs.Range("F3").Value2 = 226
$defWs.Range("G3").Value2 = 75
$defWs.Range("H3").Value2 = 52
$defWs.Range("I3").Value2 = 122
$defWs.Range("J3").Value2 = 92
$defWs.Range("N3").Value2 = 40

# ---------------------------------------------------------------------------
# ST sheet: updated running totals, plus appended samples to the running
# per-kick lists.
# ---------------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value2 = 157
$stWs.Range("D2").Value2 = 103
$stWs.Range("F2").Value2 = 62
$stWs.Range("G2").Value2 = 60
$stWs.Range("L2").Value2 = 20
$stWs.Range("M2").Value2 = 10

$stWs.Range("B3").Value2 = 79

$stWs.Range("B4").Value2 = $stWs.Range("B4").Value2 + " 67 55 65"
$stWs.Range("B5").Value2 = $stWs.Range("B5").Value2 + " 19 12 23"
$stWs.Range("B6").Value2 = $stWs.Range("B6").Value2 + " 21 0 13 6"
$stWs.Range("D3").Value2 = $stWs.Range("D3").Value2 + " 38 35 39 41 43"
$stWs.Range("D4").Value2 = $stWs.Range("D4").Value2 + " 0 0 0 0 0"
$stWs.Range("D5").Value2 = $stWs.Range("D5").Value2 + " 8 1 6"

# ---------------------------------------------------------------------------
# TURNS sheet: updated running totals.
# ---------------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B2").Value2 = 8
$turnsWs.Range("E2").Value2 = 18
$turnsWs.Range("D3").Value2 = 13

# ---------------------------------------------------------------------------
# PEN sheet: updated running totals.
# ---------------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B3").Value2 = 37
